$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.034.44"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.678.29"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.53"
$ws.Range("E5").Value = "  -4.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.77"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("E9").Value = "  -4.90%  "
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -5.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.12"
$ws.Range("E12").Value = "  -12.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.149.48"
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.05"
$ws.Range("E14").Value = "  -4.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.886.30"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("E16").Value = "  -5.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.679.05"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  -6.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.83"
$ws.Range("E20").Value = "  -5.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("E21").Value = "  -6.03%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.505"
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.33"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.12"
$ws.Range("E27").Value = "  -5.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0854"
$ws.Range("E28").Value = "  -7.06%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.34"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.00"
$ws.Range("E31").Value = "  -5.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.18"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.48"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "339.03"
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.933"
$ws.Range("E39").Value = "  -7.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.11"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.05"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("E42").Value = "  -7.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.32"
$ws.Range("E43").Value = "  -6.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.68"
$ws.Range("E44").Value = "  -8.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0560"
$ws.Range("E45").Value = "  -6.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.615"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0970"
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.72"
$ws.Range("E50").Value = "  -6.34%  "
$ws.Range("E51").Value = "  -6.08%  "
